$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level bookViews window position / active tab
# ---------------------------------------------------------------------------
# activeTab 4 ("Staircase", 0-based) -> 7 ("License (GNU LGPL)", 0-based)
# The active/selected sheet becomes the last one (License).

# ---------------------------------------------------------------------------
# 2) Rename "Licence (GNU LGPL)" -> "License (GNU LGPL)"
# ---------------------------------------------------------------------------
$wsLicense = $wb.Worksheets.Item(8)
$wsLicense.Name = "License (GNU LGPL)"

# ---------------------------------------------------------------------------
# 3) Sheet1 "Inventory - Manufacturing": set a bunch of uncertainty rows
#    (columns J..M) back to "unknown" (J=0, K/L/M="(Unknown)"), matching the
#    style already used by sibling rows (style index 10 fill/border).
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item(1)

$rows = @(19,20,27,28,49,50,62,63,64,83,84,85,86,88,89,92,112,115,129,150,156)
foreach ($r in $rows) {
    $wsInv.Cells.Item($r, 10).Value = 0
    $wsInv.Cells.Item($r, 11).Value = "(Unknown)"
    $wsInv.Cells.Item($r, 12).Value = "(Unknown)"
    $wsInv.Cells.Item($r, 13).Value = "(Unknown)"
}

# Apply the shared "unknown uncertainty" style (index 10 in this workbook's
# cellXfs) to J:M for each of those rows, matching neighbouring rows 14-16.
foreach ($r in $rows) {
    $wsInv.Range("J" + $r + ":M" + $r).Style = $wsInv.Range("J14").Style
}

# Selection moved from C8 to F26 on this sheet.
$wsInv.Range("F26").Select()

# ---------------------------------------------------------------------------
# 4) "Staircase" sheet (5th tab) is no longer the tab shown as selected.
# ---------------------------------------------------------------------------
$wsStair = $wb.Worksheets.Item(5)
$wsStair.Select()

# ---------------------------------------------------------------------------
# 5) "License (GNU LGPL)" sheet view: becomes the active/selected tab, with
#    a different selected cell (I190) while keeping the same topLeftCell
#    scroll position (A133).
# ---------------------------------------------------------------------------
$wsLicense.Activate()
$wsLicense.Range("I190").Select()
